# Step 4 results ("ÉTAPE 4" / column G) added to each of the three analysis
# sheets, copying the number format from the preceding "ÉTAPE 3" column (F)
# so the new figures render identically to the existing ones.

$wb = $excel.ActiveWorkbook

# --- Sheet "LightHouse - Portable" ---------------------------------------
$ws1 = $wb.Worksheets.Item("LightHouse - Portable")
$ws1.Range("F4").Copy($ws1.Range("G4"))
$ws1.Range("G4").Value = 81
$ws1.Range("F5").Copy($ws1.Range("G5"))
$ws1.Range("G5").Value = 88
$ws1.Range("F6").Copy($ws1.Range("G6"))
$ws1.Range("G6").Value = 87
$ws1.Range("F7").Copy($ws1.Range("G7"))
$ws1.Range("G7").Value = 78
$ws1.Range("G4:G7").Select()

# --- Sheet "LightHouse - Bureau" -----------------------------------------
$ws2 = $wb.Worksheets.Item("LightHouse - Bureau")
$ws2.Range("F4").Copy($ws2.Range("G4"))
$ws2.Range("G4").Value = 93
$ws2.Range("F5").Copy($ws2.Range("G5"))
$ws2.Range("G5").Value = 86
$ws2.Range("F6").Copy($ws2.Range("G6"))
$ws2.Range("G6").Value = 93
$ws2.Range("F7").Copy($ws2.Range("G7"))
$ws2.Range("G7").Value = 90
$ws2.Range("G4,G6,G7").Select()

# --- Sheet "GTmetrix - Bureau" -------------------------------------------
$ws3 = $wb.Worksheets.Item("GTmetrix - Bureau")
$ws3.Range("F4").Copy($ws3.Range("G4"))
$ws3.Range("G4").Value = 98
$ws3.Range("F5").Copy($ws3.Range("G5"))
$ws3.Range("G5").Value = 95
$ws3.Range("G4:G5").Select()
